$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(6).ColumnWidth = 25
Write-Host ($ws.Columns.Item(6).ColumnWidth())
